# Applies the LOM3248 course-catalog update:
#  - Updates the "Ativacao" date from 01/01/2012 to 01/01/2023
#    (this date is duplicated, by a pre-existing data quirk, into the
#    B/C cells of the "Programa resumido:" row as well)
#  - Adds English-language descriptive text into the B/C columns of the
#    "Objectives:", "Short syllabus:" and "Syllabus:" rows, which
#    previously had no content in those columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 8 / Row 13: Ativacao date 01/01/2012 -> 01/01/2023 ----
# Force Text format before assignment so Excel stores the literal text
# "01/01/2023" instead of auto-converting it to a date serial number,
# then restore the original look (font/alignment/number format) of the
# cell by copy-pasting formats from an already-correctly-styled cell.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2023"
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "01/01/2023"
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "01/01/2023"
$ws.Range("B7").Copy()
$ws.Range("B13").PasteSpecial(-4122)

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "01/01/2023"
$ws.Range("C7").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# ---- Row 11: Objectives: add English description ----
$objectivesText = "Complement students' training by addressing, in greater depth, current and relevant topics and updating with state-of-the-art topics."
$ws.Range("B11").Value = $objectivesText
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)

$ws.Range("C11").Value = $objectivesText
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)

# ---- Row 14: Short syllabus: add English description ----
$shortSyllabusText = "To be defined, according to the programmed topic."
$ws.Range("B14").Value = $shortSyllabusText
$ws.Range("B10").Copy()
$ws.Range("B14").PasteSpecial(-4122)

$ws.Range("C14").Value = $shortSyllabusText
$ws.Range("C10").Copy()
$ws.Range("C14").PasteSpecial(-4122)

# ---- Row 16: Syllabus: add English description ----
$syllabusText = "The content of this elective course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."
$ws.Range("B16").Value = $syllabusText
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)

$ws.Range("C16").Value = $syllabusText
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
